$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the string-valued (cluster) cells first, while rows 6:9 are
# still present, so that every shared string stays referenced at all times
# (this avoids the shared-strings table being compacted/reordered when the
# now-unused rows are removed further down).
$ws.Range("D2").Value = "Resolving-Mac"

$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "ECs"

$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"

$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"

# --- Row 2 numeric updates
$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 0.009727666666666667
$ws.Range("N2").Value = 0.029183
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.959538670080445
$ws.Range("R2").Value = 17.635848030724
$ws.Range("S2").Value = 0.4833500233086392
$ws.Range("T2").Value = 0.4833500233086393

# --- Row 3 numeric updates
$ws.Range("G3").Value = 65.41736466666667
$ws.Range("H3").Value = 196.252094
$ws.Range("I3").Value = 0.1569674599353791
$ws.Range("J3").Value = 0.1569674599353792
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.009727666666666667
$ws.Range("N3").Value = 0.029183
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.6363583176891112
$ws.Range("R3").Value = 5.727224859202
$ws.Range("S3").Value = 0.1569674599353791
$ws.Range("T3").Value = 0.1569674599353792

# --- Row 4 numeric updates
$ws.Range("G4").Value = 60.43484133333334
$ws.Range("H4").Value = 181.304524
$ws.Range("I4").Value = 0.1450120099461104
$ws.Range("J4").Value = 0.1450120099461104
$ws.Range("M4").Value = 0.009727666666666667
$ws.Range("N4").Value = 0.029183
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.5878899915435557
$ws.Range("R4").Value = 5.291009923892001
$ws.Range("S4").Value = 0.1450120099461104
$ws.Range("T4").Value = 0.1450120099461104

# --- Row 5 numeric updates
$ws.Range("G5").Value = 89.46554166666668
$ws.Range("H5").Value = 268.396625
$ws.Range("I5").Value = 0.2146705068098712
$ws.Range("J5").Value = 0.2146705068098712
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009727666666666667
$ws.Range("N5").Value = 0.029183
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.8702909674861113
$ws.Range("R5").Value = 7.832618707375001
$ws.Range("S5").Value = 0.2146705068098712
$ws.Range("T5").Value = 0.2146705068098712

# --- Finally remove the now-obsolete rows 6:9 (the TPM update shrank the
# data set from 8 rows to 4 rows), which also shrinks the sheet dimension
# from A1:T9 down to A1:T5.
$ws.Range("A6:T9").Delete()
